$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to append: 5 new rows (rows 42-46) continuing the daily tracker log
# for date serial 45915 (2025-09-15), following the same pattern as the
# preceding day blocks (G2-G6 goals, Progress decaying by ~1%/day,
# Percentage stays 0, Change is -0.01).
$goals = @(
    @{ Id = "G2"; Name = "Workout" },
    @{ Id = "G3"; Name = "Eat Healthy" },
    @{ Id = "G4"; Name = "Read Book" },
    @{ Id = "G5"; Name = "Investment Plan" },
    @{ Id = "G6"; Name = "Spend 10 Hours without phone" }
)

$dateSerial = 45915
$progress = 0.9234832224823122
$percentage = 0
$change = -0.01

$startRow = 42
for ($i = 0; $i -lt $goals.Count; $i++) {
    $r = $startRow + $i
    $goal = $goals[$i]

    $ws.Cells.Item($r, 1).Value = $goal.Id
    $ws.Cells.Item($r, 2).Value = $goal.Name

    $cDate = $ws.Cells.Item($r, 3)
    $cDate.Value = $dateSerial
    $cDate.NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value = $progress
    $ws.Cells.Item($r, 5).Value = $percentage
    $ws.Cells.Item($r, 6).Value = $change
}
